$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case row (row 16), mirroring the format of the existing rows.
# Copy the date cell's formatting from the row above so the new date cell
# reuses the existing "short date" style instead of creating a new one.
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A16").Value = 11
$ws.Range("B16").Value = 44956
$ws.Range("C16").Value = "Datos que superen las condiciones del form"
$ws.Range("D16").Value = "El programa deberia alojar un error de superacion de caracteres."
$ws.Range("E16").Value = "El error es mostrado"
$ws.Range("F16").Value = "Aprobado"

# Restore the selection, shifted down one row same as it was before.
$null = $ws.Range("D18").Select()
